$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 154, pushing existing rows 154:190 down to 155:191,
# and carrying row 153's formatting (D column date style) into the new row.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with this week's record. The rest of
# the columns (A, B, C, E, F, G, H, I, O, R) are constant for this market /
# product subset, matching the rows immediately above and below.
$ws.Cells.Item(154, 1).Value = 7
$ws.Cells.Item(154, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value = "Ñuble"
$ws.Cells.Item(154, 4).Value = 44543
$ws.Cells.Item(154, 5).Value = 16
$ws.Cells.Item(154, 6).Value = 100112009
$ws.Cells.Item(154, 7).Value = "Acelga"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 100
$ws.Cells.Item(154, 11).Value = 300
$ws.Cells.Item(154, 12).Value = 350
$ws.Cells.Item(154, 13).Value = 325
$ws.Cells.Item(154, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(154, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(154, 16).Value = 325
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
